# Metodo de Lopez: quitar la fila "Z&N" de cada grupo de controlador (PI y PID)
# y actualizar los valores de "b" de coma decimal.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Elimina la fila Z&N del grupo PI (fila 8) y la fila Z&N del grupo PID
# (fila 12 original, que pasa a ser la fila 11 tras el primer borrado).
$ws.Rows("8").Delete()
$ws.Rows("11").Delete()

# Usa coma decimal en los valores de "b" para PI (IAE/ISE/ITAE)
$ws.Range("D5").Value = "0,986"
$ws.Range("D6").Value = "0,952"
$ws.Range("D7").Value = "0,917"

$ws.Range("C2:D4").Select() | Out-Null
